$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 84, shifting existing rows 84:191 down to 85:192
$ws.Rows("84:84").Insert()

# Populate the newly inserted row 84 with the new weekly data point
$ws.Cells.Item(84, 1).Value = 10
$ws.Cells.Item(84, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(84, 3).Value = "La Araucanía"
$ws.Cells.Item(84, 4).Value = 45128
$ws.Cells.Item(84, 4).NumberFormat = $ws.Cells.Item(85, 4).NumberFormat
$ws.Cells.Item(84, 5).Value = 9
$ws.Cells.Item(84, 6).Value = 100112035
$ws.Cells.Item(84, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(84, 8).Value = "Sin especificar"
$ws.Cells.Item(84, 9).Value = "Primera"
$ws.Cells.Item(84, 10).Value = 60
$ws.Cells.Item(84, 11).Value = 25000
$ws.Cells.Item(84, 12).Value = 25000
$ws.Cells.Item(84, 13).Value = 25000
$ws.Cells.Item(84, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(84, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(84, 16).Value = 1667
$ws.Cells.Item(84, 17).Value = 15
$ws.Cells.Item(84, 18).Value = "Hortaliza"
